$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Title VI –" -> "Title VI" for every appropriations row in the Title VI
# block (column A, rows 14-26).
$ws.Range("A14:A26").Value = "Title VI"

# Leave the selection on the edited range, matching the saved sheet view.
[void]$ws.Range("A14:A26").Select()
